$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row labels: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$headers = @(
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
  "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
  "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Turn the data range A1:U70 into an Excel Table ("Table1") with an AutoFilter
$rng = $ws.Range("A1:U70")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# 3. Freeze the header row (split under row 1)
$ws.Activate()
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
